$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing last row (row 310) E/F values ---
$ws.Range("E310").Value = 10.027
$ws.Range("F310").Value = 10.057

# --- Append new rows 311-313 with data, copying formatting from row 310 ---

# Row 311
$ws.Range("A310:G310").Copy()
$ws.Range("A311:G311").PasteSpecial(-4122)
$ws.Range("A311").Value = 45047.33333333334
$ws.Range("B311").Value = "FX_IDC:USDMAD"
$ws.Range("C311").Value = 10.0025
$ws.Range("D311").Value = 10.2563
$ws.Range("E311").Value = 9.968400000000001
$ws.Range("F311").Value = 10.1758
$ws.Range("G311").Value = 0

# Row 312
$ws.Range("A310:G310").Copy()
$ws.Range("A312:G312").PasteSpecial(-4122)
$ws.Range("A312").Value = 45078.33333333334
$ws.Range("B312").Value = "FX_IDC:USDMAD"
$ws.Range("C312").Value = 10.1628
$ws.Range("D312").Value = 10.225
$ws.Range("E312").Value = 9.850199999999999
$ws.Range("F312").Value = 9.8803
$ws.Range("G312").Value = 0

# Row 313
$ws.Range("A310:G310").Copy()
$ws.Range("A313:G313").PasteSpecial(-4122)
$ws.Range("A313").Value = 45110.33333333334
$ws.Range("B313").Value = "FX_IDC:USDMAD"
$ws.Range("C313").Value = 9.8973
$ws.Range("D313").Value = 9.912800000000001
$ws.Range("E313").Value = 9.7285
$ws.Range("F313").Value = 9.778499999999999
$ws.Range("G313").Value = 0

Write-Host ("New dimension: " + $ws.UsedRange.Address())
